# Add a new data row (row 3) to Sheet1, mirroring the existing table layout
# (FirstName, LastName, Password, Email), and make sure the sheet view is
# left-to-right (matches the author's explicit rightToLeft="0" setting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Matthew"
$ws.Range("B3").Value = "Davis"
$ws.Range("C3").Value = "Test@1234"
$ws.Range("D3").Value = "MatthewDavis@yopmail.com"

# Ensure the sheet is displayed left-to-right, as in the target workbook.
$excel.ActiveWindow.DisplayRightToLeft = $false
